# aggiornamento fino a 1/09/2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows: row, date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
    @(358, 44432, 1, 11, 111.5392415331576),
    @(359, 44433, 1, 12, 121.6791725816264),
    @(360, 44434, 1, 12, 121.6791725816264),
    @(361, 44435, 1, 11, 111.5392415331576),
    @(362, 44436, 0, 10, 101.3993104846887),
    @(363, 44437, 2, 10, 101.3993104846887),
    @(364, 44438, 1, 7, 70.9795173392821),
    @(365, 44439, 0, 6, 60.83958629081322),
    @(366, 44440, 0, 5, 50.69965524234435)
)

foreach ($row in $data) {
    $r = $row[0]

    # Copy the style of the cell above (column A carries the date-centric
    # border/alignment/number-format) onto the new date cell before
    # setting its value.
    $ws.Range("A" + ($r - 1)).Copy($ws.Range("A" + $r))

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
